$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: SW테스트 / 테스트 케이스 설계 수정 / 7월19일 / 7월23일 / 준비 / 2
$ws.Range("A3").Value = "SW테스트"
$ws.Range("B3").Value = "테스트 케이스 설계 수정"
$ws.Range("C3").Value = "7월19일"
$ws.Range("D3").Value = "7월23일"
$ws.Range("E3").Value = "준비"

# Row 4: SW테스트 / 테스트 케이스 수행 / 8월1일 / 8월2일 / 진행 / 3
$ws.Range("A4").Value = "SW테스트"
$ws.Range("B4").Value = "테스트 케이스 수행"
$ws.Range("C4").Value = "8월1일"
$ws.Range("D4").Value = "8월2일"
$ws.Range("E4").Value = "진행"

# Column F ("importance") stores its numbers as text in this workbook (see F2),
# so force the same text storage for the new entries without leaving a
# visible style on the cells.
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "2"
$ws.Range("F3").ClearFormats()

$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "3"
$ws.Range("F4").ClearFormats()
